$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Price-column (D) values are plain numeric-looking strings stored as
# TEXT in the source sheet (e.g. "1.00", "585.82"). Assigning such a string
# straight to .Value lets Excel's type-inference reinterpret it as a number
# (dropping the trailing zero, flipping the cell to Number type). Prefixing
# with a leading apostrophe is the standard Excel "treat as text" marker; it
# is stripped on storage and keeps the cell as literal Text, matching the
# original file. Values that already contain two dots (e.g. "67.290.40")
# can't be parsed as a number anyway, so they're left unprefixed.

$ws.Range('D2').Value = '67.290.40'
$ws.Range('E2').Value = '  +0.54%  '

$ws.Range('D3').Value = '2.498.44'
$ws.Range('E3').Value = '  +0.55%  '

$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').Value = "'585.82"
$ws.Range('E5').Value = '  +0.12%  '

$ws.Range('D6').Value = "'172.51"
$ws.Range('E6').Value = '  +2.69%  '

$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('E8').Value = '  -0.31%  '

$ws.Range('D9').Value = '2.501.47'
$ws.Range('E9').Value = '  +0.71%  '

$ws.Range('E10').Value = '  +0.76%  '

$ws.Range('E11').Value = '  +0.07%  '

$ws.Range('E12').Value = '  +0.11%  '

$ws.Range('E13').Value = '  -0.80%  '

$ws.Range('D14').Value = "'25.54"
$ws.Range('E14').Value = '  -1.31%  '

$ws.Range('D15').Value = '2.919.62'
$ws.Range('E15').Value = '  -0.04%  '

$ws.Range('D16').Value = '67.213.19'
$ws.Range('E16').Value = '  +0.51%  '

$ws.Range('E17').Value = '  -1.55%  '

$ws.Range('D18').Value = '2.496.54'
$ws.Range('E18').Value = '  +1.28%  '

$ws.Range('E19').Value = '  -4.80%  '

$ws.Range('E20').Value = '  -5.30%  '

$ws.Range('D21').Value = "'351.67"
$ws.Range('E21').Value = '  -3.02%  '

$ws.Range('E22').Value = '  -0.78%  '

$ws.Range('D23').Value = "'1.01"
$ws.Range('E23').Value = '  +0.55%  '

$ws.Range('D24').Value = "'4.25"
$ws.Range('E24').Value = '  -4.78%  '

$ws.Range('D25').Value = "'68.74"
$ws.Range('E25').Value = '  -3.01%  '

$ws.Range('E26').Value = '  -1.84%  '

$ws.Range('D27').Value = "'9.27"
$ws.Range('E27').Value = '  -2.01%  '

$ws.Range('D28').Value = "'0.999"
$ws.Range('E28').Value = '  -0.27%  '

$ws.Range('D29').Value = '2.624.52'
$ws.Range('E29').Value = '  +0.82%  '

$ws.Range('D30').Value = '0.0₃0908'
$ws.Range('E30').Value = '  -2.62%  '

$ws.Range('D31').Value = "'512.98"
$ws.Range('E31').Value = '  -0.43%  '

$ws.Range('D32').Value = "'7.86"
$ws.Range('E32').Value = '  -2.71%  '

$ws.Range('E33').Value = '  -2.45%  '

$ws.Range('E34').Value = '  -3.11%  '

$ws.Range('D35').Value = "'1.00"
$ws.Range('E35').Value = '  +0.01%  '

$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').Value = "'160.05"
$ws.Range('E36').Value = '  +1.04%  '

$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = "'0.118"
$ws.Range('E37').Value = '  -6.92%  '

$ws.Range('D38').Value = "'18.71"
$ws.Range('E38').Value = '  +0.81%  '

$ws.Range('D39').Value = "'18.29"
$ws.Range('E39').Value = '  -3.43%  '

$ws.Range('E40').Value = '  -5.35%  '

$ws.Range('E41').Value = '  -3.10%  '

$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').Value = "'1.00"
$ws.Range('E42').Value = '  -0.07%  '

$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D43').Value = "'4.85"
$ws.Range('E43').Value = '  -2.00%  '

$ws.Range('D44').Value = "'0.329"
$ws.Range('E44').Value = '  -1.55%  '

$ws.Range('D45').Value = "'2.38"
$ws.Range('E45').Value = '  -3.03%  '

$ws.Range('D46').Value = "'38.78"
$ws.Range('E46').Value = '  -1.16%  '

$ws.Range('D47').Value = "'143.55"
$ws.Range('E47').Value = '  +0.61%  '

$ws.Range('D48').Value = "'0.516"
$ws.Range('E48').Value = '  -4.10%  '

$ws.Range('D49').Value = "'3.46"
$ws.Range('E49').Value = '  -3.70%  '

$ws.Range('E50').Value = '  -5.53%  '

$ws.Range('E51').Value = '  -4.30%  '
